$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("SP24092022115848", 101, "PANADOL STRIP 20", 3, 261),
    @("SP24092022115848", 100, "PANADOL STRIP 10", 6, 540),
    @("SP24092022115848", 109, "ATORVASTATIN SANDOZ 10MG TABLETS 30", 2, 23.2),
    @("SP24092022120034", 100, "PANADOL STRIP 10", 6, 540),
    @("SP24092022120034", 101, "PANADOL STRIP 20", 3, 261),
    @("SP24092022120034", 180, "EPIDUO GEL 30g PUMP", 2, 110.76),
    @("SP24092022160310", 100, "PANADOL STRIP 10", 4, 240),
    @("SP24092022160310", 102, "NUROFEN STRIP 15", 1, 12),
    @("SP24092022160509", 102, "NUROFEN STRIP 15", 1, 12)
)

$startRow = 24
$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
